$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JUROS E MULTA")

# --- Widen column D to fit the new, longer descriptions ---
$ws.Columns("D").ColumnWidth = 50.8

# --- Row 16: Total row (write the "Total" label first so it gets the lowest new
#     shared-string index, matching how the sheet was actually authored) ---
$ws.Range("D10").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E10").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("D16").Value = "Total"

# --- Row 14: IPTU proportional value ---
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E9").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("D14").Value = "IPTU 2023 (R`$ 926,00) proporcional"
$ws.Range("E14").Formula = "=926/12*4"

# --- Row 15: Condominio proportional value ---
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E9").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("D15").Value = "Condomínio (25/08 até  31/08) (R`$610,54) proporcional"
$ws.Range("E15").Formula = "=610.54/30*6"

# --- Total formula, now that E14/E15 exist ---
$ws.Range("E16").Formula = "=SUM(E14:E15)"

# --- Give the new total row its own (gray) fill color ---
$ws.Range("D16:E16").Interior.ThemeColor = 2
$ws.Range("D16:E16").Interior.TintAndShade = -0.249977111117893

# --- Recalculate everything so cached formula results are correct ---
$excel.CalculateFull()

# --- Update the active selection to match the end of the editing session ---
$ws.Range("D53").Select()

$wb.Save()
